# "Add results for Sioux Falls 2017." — refresh the Column D result figures
# on the (only) worksheet, and nudge the saved scroll position down a few
# rows (topLeftCell A25 -> A28) to match where the user left the view.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New Column D values (Column E / labels / everything else is unchanged).
$updates = @{
    "D3"  = 13962
    "D4"  = 10035
    "D5"  = 8545
    "D6"  = 140
    "D8"  = 3974
    "D9"  = 3049
    "D10" = 837
    "D11" = 446
    "D13" = 2881
    "D14" = 2862
    "D15" = 1294
    "D16" = 221
    "D19" = 112
    "D22" = 1833
    "D23" = 1968
    "D24" = 243
    "D26" = 4743
    "D28" = 539
    "D30" = 3642
    "D31" = 382
    "D33" = 851
    "D34" = 67
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

# Best-effort: move the window's scroll anchor to row 28 (was row 25) without
# touching the current selection (stays D35, per the saved view state).
try {
    $win = $wb.Windows.Item(1)
    $win.ScrollRow = 28
    $win.ScrollColumn = 1
} catch {
    # Scroll position is cosmetic only; ignore if unsupported.
}
